$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected values for kinetic model ddsAAR26 (row 18)
$ws.Range("E18").Value = 573
$ws.Range("F18").Value = 275
$ws.Range("G18").Value = 701

# Update the selected/active cell to reflect new view position
$ws.Range("E19").Select()
